# Updated cryptos list on Sat Jul 22 21:53:11 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) snapshot values
# for the coinranking.com crypto table, and fixes the WrappedEther/Polygon
# row ordering (rows 12-13 had swapped places).
#
# NOTE: Several Price values look like plain decimals (e.g. "242.42",
# "1.000") which Excel's literal-entry parser would otherwise interpret as
# numbers. The source data stores these as literal text (some prices even
# use two "thousands" dots, e.g. "29.836.83", which is why the whole column
# is text). A leading apostrophe forces those ambiguous values to stay text,
# matching the original cell type/content - it is omitted wherever the
# string already can't be parsed as a number (so it is stored as text with
# no extra styling, same as before).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.836.83"
$ws.Range("E2").Value = "  -0.33%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.887.77"
$ws.Range("E3").Value = "  -0.35%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - XRP
$ws.Range("D5").Value = "'0.7524"
$ws.Range("E5").Value = "  -2.64%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'242.42"

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.12%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3125"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 - Solana
$ws.Range("D9").Value = "'25.32"
$ws.Range("E9").Value = "  -1.32%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.07115"
$ws.Range("E10").Value = "  -2.98%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.08483"
$ws.Range("E11").Value = "  +4.99%  "

# Row 12 - was WrappedEther, now Polygon
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7600"
$ws.Range("E12").Value = "  -1.45%  "

# Row 13 - was Polygon, now WrappedEther
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.897.83"
$ws.Range("E13").Value = "  -0.68%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -2.10%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "'93.36"
$ws.Range("E15").Value = "  -0.67%  "

# Row 16 - Uniswap
$ws.Range("D16").Value = "'6.129"
$ws.Range("E16").Value = "  -1.26%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "29.974.30"
$ws.Range("E17").Value = "  +0.19%  "

# Row 18 - Avalanche
$ws.Range("E18").Value = "  -1.64%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'243.40"
$ws.Range("E19").Value = "  -1.42%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "'0.000007807"
$ws.Range("E20").Value = "  -0.02%  "

# Row 21 - Dai
$ws.Range("D21").Value = "'0.9988"
$ws.Range("E21").Value = "  -0.15%  "

# Row 22 - WrappedliquidstakedEther2.0
$ws.Range("D22").Value = "2.139.17"
$ws.Range("E22").Value = "  -0.79%  "

# Row 23 - Chainlink
$ws.Range("D23").Value = "'8.029"
$ws.Range("E23").Value = "  -0.89%  "

# Row 24 - BinanceUSD
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.11%  "

# Row 25 - Stellar
$ws.Range("D25").Value = "'0.1593"
$ws.Range("E25").Value = "  +1.17%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "'9.371"
$ws.Range("E26").Value = "  -0.79%  "

# Row 27 - Monero
$ws.Range("D27").Value = "'162.70"

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +0.18%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "'2.031"
$ws.Range("E29").Value = "  +0.29%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  +3.50%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "'1.536"
$ws.Range("E31").Value = "  -0.60%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'4.504"
$ws.Range("E32").Value = "  +0.82%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "'4.160"
$ws.Range("E33").Value = "  +2.51%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "'0.05428"
$ws.Range("E34").Value = "  -2.19%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  +0.28%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "'0.7519"
$ws.Range("E36").Value = "  +0.14%  "

# Row 37 - Frax
$ws.Range("D37").Value = "'1.004"
$ws.Range("E37").Value = "  +0.17%  "

# Row 38 - HuobiToken
$ws.Range("D38").Value = "'2.711"

# Row 39 - VeChain
$ws.Range("D39").Value = "'0.01944"
$ws.Range("E39").Value = "  +0.79%  "

# Row 40 - MXToken
$ws.Range("D40").Value = "'2.771"
$ws.Range("E40").Value = "  -0.66%  "

# Row 41 - TheSandbox
$ws.Range("D41").Value = "'0.4462"
$ws.Range("E41").Value = "  -0.02%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "'6.109"
$ws.Range("E42").Value = "  +2.45%  "

# Row 43 - Maker
$ws.Range("D43").Value = "1.097.73"
$ws.Range("E43").Value = "  -0.30%  "

# Row 44 - Aave
$ws.Range("D44").Value = "'72.71"
$ws.Range("E44").Value = "  -1.63%  "

# Row 45 - TrustWalletToken
$ws.Range("D45").Value = "'0.8603"
$ws.Range("E45").Value = "  +1.15%  "

# Row 46 - PaxDollar
$ws.Range("D46").Value = "'0.9999"
$ws.Range("E46").Value = "  -0.13%  "

# Row 47 - Aptos
$ws.Range("D47").Value = "'7.720"
$ws.Range("E47").Value = "  +2.78%  "

# Row 48 - Quant
$ws.Range("D48").Value = "'102.59"
$ws.Range("E48").Value = "  +0.38%  "

# Row 49 - RenderToken
$ws.Range("D49").Value = "'1.860"
$ws.Range("E49").Value = "  -1.27%  "

# Row 50 - SynthetixNetwork
$ws.Range("D50").Value = "'3.043"
$ws.Range("E50").Value = "  +1.68%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.039.14"
$ws.Range("E51").Value = "  +0.09%  "
